$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add SUM formula totals row (row 6) under the existing data (rows 2-5)
$ws.Range("B6").Formula = "=SUM(B2:B5)"
$ws.Range("C6").Formula = "=SUM(C2:C5)"
$ws.Range("D6").Formula = "=SUM(D2:D5)"
$ws.Range("E6").Formula = "=SUM(E2:E5)"
$ws.Range("F6").Formula = "=SUM(F2:F5)"
$ws.Range("G6").Formula = "=SUM(G2:G5)"
$ws.Range("H6").Formula = "=SUM(H2:H5)"

# Format the new totals row with the Currency cell style
$ws.Range("B6:H6").Style = "Currency"

# Clear the chart's explicit style reference (use workbook default style)
$chart = $ws.ChartObjects().Item(1).Chart
$chart.ChartStyle = $null
